# Michaels stress inversion technique added
# Applies the data/formatting changes to Sheet1 of the Angelier fault-slip
# dataset: a batch of individual measurement-value corrections, a new
# "TYM" site record inserted among rows 69-72 (shifting the three
# following records down and dropping the old trailing one), a bottom
# border added under the "TYM"/row-34 summary rows, and the active
# selection/scroll position updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Individual cell value corrections -------------------------------
$ws.Range("F5").Value  = 84
$ws.Range("D17").Value = 85
$ws.Range("D20").Value = 58
$ws.Range("F21").Value = 78
$ws.Range("D34").Value = 80
$ws.Range("F36").Value = 80
$ws.Range("D37").Value = 80
$ws.Range("F38").Value = 78
$ws.Range("D40").Value = 88
$ws.Range("D41").Value = 78
$ws.Range("F44").Value = 88
$ws.Range("F54").Value = 69
$ws.Range("F58").Value = 89
$ws.Range("D64").Value = 58
$ws.Range("F66").Value = 84
$ws.Range("D67").Value = 68
$ws.Range("F67").Value = 87

# --- New TYM record inserted at row 69, pushing the former rows
#     69-71 down to 70-72 (with minor pitch corrections) and
#     dropping the former row 72 entirely. -----------------------------
$ws.Range("A69").Value = "TYM"
$ws.Range("B69").Value = "N"
$ws.Range("C69").Value = 70
$ws.Range("D69").Value = 69
$ws.Range("E69").Value = "N"
$ws.Range("F69").Value = 88
$ws.Range("G69").Value = "W"
$ws.Range("H69:K69").ClearContents()

$ws.Range("A70").Value = "TYM"
$ws.Range("B70").Value = "N"
$ws.Range("C70").Value = 17
$ws.Range("D70").Value = 68
$ws.Range("E70").Value = "W"
$ws.Range("F70").Value = 80
$ws.Range("G70").Value = "N"
$ws.Range("H70").Value = 37
$ws.Range("I70").Value = 14
$ws.Range("J70").Value = 48
$ws.Range("K70").Value = 0

$ws.Range("A71").Value = "TYM"
$ws.Range("B71").Value = "N"
$ws.Range("C71").Value = 96
$ws.Range("D71").Value = 70
$ws.Range("E71").Value = "N"
$ws.Range("F71").Value = 72
$ws.Range("G71").Value = "W"
$ws.Range("H71").Value = 50
$ws.Range("I71").Value = 6
$ws.Range("J71").Value = 57
$ws.Range("K71").Value = 3

$ws.Range("A72").Value = "TYM"
$ws.Range("B72").Value = "N"
$ws.Range("C72").Value = 89
$ws.Range("D72").Value = 71
$ws.Range("E72").Value = "N"
$ws.Range("F72").Value = 89
$ws.Range("G72").Value = "E"
$ws.Range("H72").Value = 29
$ws.Range("I72").Value = 5
$ws.Range("J72").Value = 43
$ws.Range("K72").Value = 4

# --- Bottom border under the "AVB" (row 34) and new "TYM" (row 72)
#     site-group summary rows ------------------------------------------
$ws.Range("A34:K34").Borders.Item(9).LineStyle = 1
$ws.Range("A34:K34").Borders.Item(9).Weight = 2
$ws.Range("A72:K72").Borders.Item(9).LineStyle = 1
$ws.Range("A72:K72").Borders.Item(9).Weight = 2

# --- Selection / scroll position update -------------------------------
$ws.Range("E1").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
